$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Source fixtures" / "Source type" No answers for the first
# (default) calibration column, mirroring the existing "Yes" entries in
# column B for rows 4 and 6.
$ws.Range("B3").Value = "No"
$ws.Range("B5").Value = "No"

# Update the saved selection/active cell on the sheet.
$ws.Range("C3").Select()
